# Adds a new "cost-fix" column to both the "Commodity" and "Process" sheets,
# matching the upstream rivus commit that introduced fixed (size-independent)
# maintenance costs alongside the existing variable costs.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Commodity sheet: insert new column E "cost-fix" (between cost-inv-var
# and cost-var), values are plain integers using the surrounding style.
# ---------------------------------------------------------------------
$wsCommodity = $wb.Worksheets.Item("Commodity")

$wsCommodity.Columns.Item(5).Insert()

$wsCommodity.Range("E1").Value = "cost-fix"
$wsCommodity.Range("E2").Value = 7
$wsCommodity.Range("E3").Value = 5
$wsCommodity.Range("E4").Value = 10
$wsCommodity.Range("E5").Formula = "=NA()"

$rCommodity = $wsCommodity.Range("E1")
$rCommodity.Validation.Add(0)
$rCommodity.Validation.InputTitle = "Variable fixed costs (€/kW/m)"
$rCommodity.Validation.InputMessage = "Capacity-dependent fixed costs for maintaining transmission capacity."
$rCommodity.Validation.ShowInput = $true
$rCommodity.Validation.ShowError = $true

# ---------------------------------------------------------------------
# Process sheet: insert new column D "cost-fix" (between cost-inv-var and
# cost-var), formatted with the "0.0" number format (new style, fill
# matching the header/data convention used by the other numeric columns).
# Also fixes the "District heating plant" cost-inv-var value (2 -> 600).
# ---------------------------------------------------------------------
$wsProcess = $wb.Worksheets.Item("Process")

$wsProcess.Columns.Item(4).Insert()

$wsProcess.Range("D1").Value = "cost-fix"
$wsProcess.Range("D1").NumberFormat = "0.0"
$wsProcess.Range("D2:D8").NumberFormat = "0.0"
$wsProcess.Range("D2").Value = 50
$wsProcess.Range("D3").Value = 30
$wsProcess.Range("D4").Value = 75
$wsProcess.Range("D5").Value = 100
$wsProcess.Range("D6").Value = 60
$wsProcess.Range("D7").Value = 50
$wsProcess.Range("D8").Value = 80

# District heating plant: cost-inv-var corrected from 2 to 600
$wsProcess.Range("C6").Value = 600

$rProcess = $wsProcess.Range("D1")
$rProcess.Validation.Add(0)
$rProcess.Validation.InputTitle = "Specific fixed costs (€/kW)"
$rProcess.Validation.InputMessage = "Size-dependent part for maintaining a plant."
$rProcess.Validation.ShowInput = $true
$rProcess.Validation.ShowError = $true

# ---------------------------------------------------------------------
# Selection / active sheet bookkeeping, matching the saved workbook state:
# the active sheet moves from "Process" to "Commodity", and each sheet
# remembers a new last-used selection.
# ---------------------------------------------------------------------
$wsProcess.Range("E8").Select()

$wsCommodity.Activate()
$wsCommodity.Range("H25").Select()
